$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data cells for rows 2-5 with new sensor readings (custom accuracy refresh) ---
$ws.Range("A2").Value = 45171.50694444445
$ws.Range("B2").Value = 14.835
$ws.Range("C2").Value = 9.791
$ws.Range("D2").Value = 3.698
$ws.Range("E2").Value = 32.243
$ws.Range("F2").Value = 24.166
$ws.Range("G2").Value = 11.51
$ws.Range("H2").Value = 34.958
$ws.Range("I2").Value = 18.033
$ws.Range("J2").Value = 7.29
$ws.Range("K2").Value = 10.735
$ws.Range("L2").Value = 12.533
$ws.Range("M2").Value = 13.25
$ws.Range("N2").Value = 3.739
$ws.Range("O2").Value = 11.655
$ws.Range("P2").Value = 16.06
$ws.Range("Q2").Value = 10.282
$ws.Range("R2").Value = 3.096
$ws.Range("S2").Value = 1.74
$ws.Range("T2").Value = 170.025
$ws.Range("U2").Value = 32.298
$ws.Range("V2").Value = 10.758
$ws.Range("W2").Value = 20.812
$ws.Range("X2").Value = 10.713
$ws.Range("Y2").Value = 2.837
$ws.Range("Z2").Value = 18.288
$ws.Range("AA2").Value = 9.502000000000001
$ws.Range("AB2").Value = 8.641999999999999
$ws.Range("AC2").Value = 10.303
$ws.Range("AD2").Value = 12.679
$ws.Range("AE2").Value = 3.311
$ws.Range("AF2").Value = 31.418
$ws.Range("AG2").Value = 5.68
$ws.Range("AH2").Value = 13.449

$ws.Range("A3").Value = 45171.51388888889
$ws.Range("B3").Value = 23.968
$ws.Range("C3").Value = 17.405
$ws.Range("D3").Value = 2.03
$ws.Range("E3").Value = 52.423
$ws.Range("F3").Value = 41.972
$ws.Range("G3").Value = 18.744
$ws.Range("H3").Value = 70.51300000000001
$ws.Range("I3").Value = 29.086
$ws.Range("J3").Value = 12.717
$ws.Range("K3").Value = 18.715
$ws.Range("L3").Value = 20.86
$ws.Range("M3").Value = 22.106
$ws.Range("N3").Value = 6.038
$ws.Range("O3").Value = 18.798
$ws.Range("P3").Value = 26.594
$ws.Range("Q3").Value = 16.059
$ws.Range("R3").Value = 1.528
$ws.Range("S3").Value = 1.248
$ws.Range("T3").Value = 278.816
$ws.Range("U3").Value = 52.546
$ws.Range("V3").Value = 17.351
$ws.Range("W3").Value = 35.016
$ws.Range("X3").Value = 18.321
$ws.Range("Y3").Value = 3.09
$ws.Range("Z3").Value = 34.931
$ws.Range("AA3").Value = 15.326
$ws.Range("AB3").Value = 13.686
$ws.Range("AC3").Value = 16.132
$ws.Range("AD3").Value = 21.641
$ws.Range("AE3").Value = 1.246
$ws.Range("AF3").Value = 64.446
$ws.Range("AG3").Value = 9.625999999999999
$ws.Range("AH3").Value = 21.693

$ws.Range("A4").Value = 45171.52083333334
$ws.Range("B4").Value = 10.04
$ws.Range("C4").Value = 7.173
$ws.Range("D4").Value = 1.096
$ws.Range("E4").Value = 22.099
$ws.Range("F4").Value = 17.318
$ws.Range("G4").Value = 7.811
$ws.Range("H4").Value = 34.781
$ws.Range("I4").Value = 12.216
$ws.Range("J4").Value = 5.283
$ws.Range("K4").Value = 7.627
$ws.Range("L4").Value = 8.763999999999999
$ws.Range("M4").Value = 9.321999999999999
$ws.Range("N4").Value = 2.539
$ws.Range("O4").Value = 7.895
$ws.Range("P4").Value = 11.128
$ws.Range("Q4").Value = 6.904
$ws.Range("R4").Value = 0.951
$ws.Range("S4").Value = 0.624
$ws.Range("T4").Value = 112.87
$ws.Range("U4").Value = 22.235
$ws.Range("V4").Value = 7.288
$ws.Range("W4").Value = 14.666
$ws.Range("X4").Value = 7.631
$ws.Range("Y4").Value = 1.454
$ws.Range("Z4").Value = 16.383
$ws.Range("AA4").Value = 6.437
$ws.Range("AB4").Value = 5.828
$ws.Range("AC4").Value = 6.859
$ws.Range("AD4").Value = 9.006
$ws.Range("AE4").Value = 0.766
$ws.Range("AF4").Value = 32.019
$ws.Range("AG4").Value = 3.965
$ws.Range("AH4").Value = 9.112

$ws.Range("A5").Value = 45171.52777777778
$ws.Range("B5").Value = 7.64
$ws.Range("C5").Value = 5.48
$ws.Range("D5").Value = 0.8100000000000001
$ws.Range("E5").Value = 16.85
$ws.Range("F5").Value = 13.2
$ws.Range("G5").Value = 5.94
$ws.Range("H5").Value = 25.24
$ws.Range("I5").Value = 9.31
$ws.Range("J5").Value = 4.02
$ws.Range("K5").Value = 5.8
$ws.Range("L5").Value = 6.69
$ws.Range("M5").Value = 7.13
$ws.Range("N5").Value = 1.93
$ws.Range("O5").Value = 6.02
$ws.Range("P5").Value = 8.460000000000001
$ws.Range("Q5").Value = 5.27
$ws.Range("R5").Value = 0.71
$ws.Range("S5").Value = 0.46
$ws.Range("T5").Value = 84.23
$ws.Range("U5").Value = 16.88
$ws.Range("V5").Value = 5.55
$ws.Range("W5").Value = 11.13
$ws.Range("X5").Value = 5.81
$ws.Range("Y5").Value = 1.11
$ws.Range("Z5").Value = 11.93
$ws.Range("AA5").Value = 4.9
$ws.Range("AB5").Value = 4.45
$ws.Range("AC5").Value = 5.23
$ws.Range("AD5").Value = 6.87
$ws.Range("AE5").Value = 0.55
$ws.Range("AF5").Value = 23.07
$ws.Range("AG5").Value = 3.02
$ws.Range("AH5").Value = 6.94

# --- Drop the old 6th timestamp row now that data has been refreshed (1000-row dataset trimmed to match) ---
$ws.Rows.Item(6).Delete() | Out-Null

# --- Widen several columns from 7 to 8 raw width units (use numeric column index; letter-keyed Columns.Item is unreliable here) ---
$ws.Columns.Item(3).ColumnWidth = 7.17
$ws.Columns.Item(7).ColumnWidth = 7.17
$ws.Columns.Item(10).ColumnWidth = 7.17
$ws.Columns.Item(11).ColumnWidth = 7.17
$ws.Columns.Item(12).ColumnWidth = 7.17
$ws.Columns.Item(13).ColumnWidth = 7.17
$ws.Columns.Item(15).ColumnWidth = 7.17
$ws.Columns.Item(16).ColumnWidth = 7.17
$ws.Columns.Item(17).ColumnWidth = 7.17
$ws.Columns.Item(22).ColumnWidth = 7.17
$ws.Columns.Item(24).ColumnWidth = 7.17
$ws.Columns.Item(27).ColumnWidth = 7.17
$ws.Columns.Item(28).ColumnWidth = 7.17
$ws.Columns.Item(29).ColumnWidth = 7.17
$ws.Columns.Item(34).ColumnWidth = 7.17
